$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (column A holds the "Beteckning" designation
# for every data row, so use it to discover the extent of the table).
$lastRow = $ws.UsedRange.Rows.Count - 1
if ($lastRow -lt 1) { $lastRow = 506 }

# 1) Update the "Förändrad" (column C) date for every data row (rows 2..lastRow)
#    from the old serial date (45184) to the new one (45186).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($val -ne $null) {
        $cell.Value = 45186
    }
}

# 2) For every HYPERLINK(...) formula that only has the URL argument, add a second
#    argument containing the row's designation (column A), turning
#    HYPERLINK("url") into HYPERLINK("url", "A 58400-2021") style links.
$cols = @("S","T","U","V","W","X","Y")
for ($r = 2; $r -le $lastRow; $r++) {
    $designation = $ws.Cells.Item($r, 1).Value()
    foreach ($colLetter in $cols) {
        $cell = $ws.Range($colLetter + $r)
        $formula = $cell.Formula
        if ($formula -and $formula.StartsWith("=HYPERLINK(")) {
            if ($formula -match '^=HYPERLINK\("([^"]*)"\)$') {
                $url = $matches[1]
                $cell.Formula = '=HYPERLINK("' + $url + '", "' + $designation + '")'
            }
        }
    }
}
